$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) of the last existing row (140) down into the
# two new rows (141-142) so the bold/bordered index column and the
# date-time number format on column E carry over correctly.
$ws.Range("A140:V140").Copy()
$ws.Range("A141:V142").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 141: FC Botosani 0 - 0 FC Rapid Bucuresti
$ws.Range("A141").Value = 140
$ws.Range("B141").Value = "romania"
$ws.Range("C141").Value = "liga-1"
$ws.Range("D141").Value = "2023-2024"
$ws.Range("E141").Value = 45262.66666666666
$ws.Range("F141").Value = "FC Botosani"
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = "FC Rapid Bucuresti"
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 4.26
$ws.Range("K141").Value = "30/11/2023 08:12"
$ws.Range("L141").Value = 4.01
$ws.Range("M141").Value = "02/12/2023 15:46"
$ws.Range("N141").Value = 3.61
$ws.Range("O141").Value = "30/11/2023 08:12"
$ws.Range("P141").Value = 3.51
$ws.Range("Q141").Value = "02/12/2023 15:39"
$ws.Range("R141").Value = 1.77
$ws.Range("S141").Value = "30/11/2023 08:12"
$ws.Range("T141").Value = 1.95
$ws.Range("U141").Value = "02/12/2023 15:46"
$ws.Range("V141").Value = "https://www.betexplorer.com/football/romania/liga-1/fc-botosani-rapid-bucuresti/ltL8Qeoc/"

# Row 142: Univ. Craiova 1 - 0 CFR Cluj
$ws.Range("A142").Value = 141
$ws.Range("B142").Value = "romania"
$ws.Range("C142").Value = "liga-1"
$ws.Range("D142").Value = "2023-2024"
$ws.Range("E142").Value = 45262.8125
$ws.Range("F142").Value = "Univ. Craiova"
$ws.Range("G142").Value = 1
$ws.Range("H142").Value = "CFR Cluj"
$ws.Range("I142").Value = 0
$ws.Range("J142").Value = 2.68
$ws.Range("K142").Value = "30/11/2023 08:12"
$ws.Range("L142").Value = 2.33
$ws.Range("M142").Value = "02/12/2023 19:22"
$ws.Range("N142").Value = 3.06
$ws.Range("O142").Value = "30/11/2023 08:12"
$ws.Range("P142").Value = 3.21
$ws.Range("Q142").Value = "02/12/2023 19:22"
$ws.Range("R142").Value = 2.68
$ws.Range("S142").Value = "30/11/2023 08:12"
$ws.Range("T142").Value = 3.29
$ws.Range("U142").Value = "02/12/2023 19:26"
$ws.Range("V142").Value = "https://www.betexplorer.com/football/romania/liga-1/univ-craiova-cfr-cluj/x2ij8EVS/"
